# Update countries & provincias Spain
# Applies the data refresh captured in the commit: a few country rows get
# revised case counts (Pakistan, Kirguistan), Israel overtakes Kazajistan
# and Islas Malvinas overtakes Montserrat in the ranking (so those two
# pairs of rows swap places), and the "last updated" timestamp moves from
# 06:09 to 07:26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 07:26"

# --- Row 19: Pakistan (row position unchanged, just refreshed stats) --
$ws.Range("B19").Value = 293711
$ws.Range("C19").Value = 450
$ws.Range("D19").Value = 278425
$ws.Range("E19").Value = 9031
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 11
$ws.Range("H19").Value = 6255

# --- Rows 32/33: Israel overtakes Kazajistan ---------------------------
# Row 32 now holds Israel (updated totals), row 33 now holds Kazajistan
# (same totals Kazajistan had before, just shifted down a row).
$ws.Range("A32").Value = "Israel"
$ws.Range("B32").Value = 105063
$ws.Range("C32").Value = 591
$ws.Range("D32").Value = 83013
$ws.Range("E32").Value = 21203
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 847

$ws.Range("A33").Value = "Kazajistan"
$ws.Range("B33").Value = 104902
$ws.Range("C33").Value = 184
$ws.Range("D33").Value = 93405
$ws.Range("E33").Value = 10082
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 1415

# --- Row 56: Kirguistan (row position unchanged, refreshed stats) -----
$ws.Range("B56").Value = 43204
$ws.Range("C56").Value = 78
$ws.Range("D56").Value = 36820
$ws.Range("E56").Value = 5327
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 1057

# --- Rows 214/215: Islas Malvinas overtakes Montserrat ------------------
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1
